# Delete the "Exercise 2" (slide 7) and "Exercise 3" (slide 8) slides from
# the deck. The trailing "[end of STA6206 BDA Practical 4]" slide (originally
# slide 9) shifts up to become the new last slide (position 7).
$p = $ppt.ActivePresentation

# Delete slide at position 8 first (Exercise 3), then position 7 (Exercise 2),
# so indices of the slide we are about to remove stay valid as the count shrinks.
$p.Slides.Item(8).Delete()
$p.Slides.Item(7).Delete()
